{"js": "// Update the header date and every \"two-digit \u00f7 one-digit\" answer cell in\n// the single table, going cell-by-cell (row/col) in document order so\n// duplicate \"before\" values (e.g. \"47\u00f74=11, 3\" appears twice) each map to\n// their own distinct replacement instead of a blanket find/replace.\n\nconst body = context.document.body;\n\n// --- 1) Header paragraph: \"2025-11-18 Tuesday\" -> \"2025-11-19 Wednesday\" ---\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst dateOld = \"2025-11-18 Tuesday\";\nconst dateNew = \"2025-11-19 Wednesday\";\nconst dateParagraph = paragraphs.items.find((p) => p.text === dateOld);\nif (dateParagraph) {\n  dateParagraph.getRange().insertText(dateNew, \"Replace\");\n}\n\n// --- 2) Table answer cells ---\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// (row, col) are 0-based; only the 5 \"data\" rows (0, 4, 8, 12, 16) of the\n// 20-row table carry text, the rest are blank spacer rows.\nconst cellUpdates = [\n  { row: 0, col: 0, oldText: \"47\u00f79=5, 2\", newText: \"92\u00f72=46, 0\" },\n  { row: 0, col: 1, oldText: \"12\u00f78=1, 4\", newText: \"23\u00f73=7, 2\" },\n  { row: 0, col: 2, oldText: \"83\u00f75=16, 3\", newText: \"98\u00f72=49, 0\" },\n  { row: 0, col: 3, oldText: \"30\u00f77=4, 2\", newText: \"33\u00f72=16, 1\" },\n  { row: 0, col: 4, oldText: \"55\u00f76=9, 1\", newText: \"23\u00f75=4, 3\" },\n  { row: 4, col: 0, oldText: \"71\u00f75=14, 1\", newText: \"48\u00f76=8, 0\" },\n  { row: 4, col: 1, oldText: \"81\u00f74=20, 1\", newText: \"50\u00f74=12, 2\" },\n  { row: 4, col: 2, oldText: \"17\u00f78=2, 1\", newText: \"97\u00f72=48, 1\" },\n  { row: 4, col: 3, oldText: \"40\u00f79=4, 4\", newText: \"91\u00f77=13, 0\" },\n  { row: 4, col: 4, oldText: \"72\u00f74=18, 0\", newText: \"90\u00f75=18, 0\" },\n  { row: 8, col: 0, oldText: \"28\u00f77=4, 0\", newText: \"81\u00f79=9, 0\" },\n  { row: 8, col: 1, oldText: \"87\u00f79=9, 6\", newText: \"65\u00f79=7, 2\" },\n  { row: 8, col: 2, oldText: \"61\u00f78=7, 5\", newText: \"89\u00f73=29, 2\" },\n  { row: 8, col: 3, oldText: \"59\u00f78=7, 3\", newText: \"73\u00f79=8, 1\" },\n  { row: 8, col: 4, oldText: \"14\u00f72=7, 0\", newText: \"86\u00f77=12, 2\" },\n  { row: 12, col: 0, oldText: \"45\u00f73=15, 0\", newText: \"15\u00f77=2, 1\" },\n  { row: 12, col: 1, oldText: \"31\u00f78=3, 7\", newText: \"91\u00f73=30, 1\" },\n  { row: 12, col: 2, oldText: \"13\u00f72=6, 1\", newText: \"43\u00f76=7, 1\" },\n  { row: 12, col: 3, oldText: \"47\u00f74=11, 3\", newText: \"93\u00f77=13, 2\" },\n  { row: 12, col: 4, oldText: \"16\u00f78=2, 0\", newText: \"64\u00f79=7, 1\" },\n  { row: 16, col: 0, oldText: \"21\u00f72=10, 1\", newText: \"56\u00f73=18, 2\" },\n  { row: 16, col: 1, oldText: \"28\u00f78=3, 4\", newText: \"44\u00f78=5, 4\" },\n  { row: 16, col: 2, oldText: \"47\u00f74=11, 3\", newText: \"65\u00f72=32, 1\" },\n  { row: 16, col: 3, oldText: \"14\u00f75=2, 4\", newText: \"31\u00f73=10, 1\" },\n  { row: 16, col: 4, oldText: \"97\u00f75=19, 2\", newText: \"14\u00f79=1, 5\" },\n];\n\nfor (const update of cellUpdates) {\n  const cell = table.getCell(update.row, update.col);\n  cell.value = update.newText;\n}\n\nawait context.sync();\n", "ps1": "# Update the header date and every \"two-digit \u00f7 one-digit\" answer cell in\n# the single table, addressing cells by (row, col) in document order so\n# duplicate \"before\" values (e.g. \"47\u00f74=11, 3\" appears twice) each map to\n# their own distinct replacement instead of a blanket find/replace.\n\n$d = $word.ActiveDocument\n\n# --- 1) Header paragraph: \"2025-11-18 Tuesday\" -> \"2025-11-19 Wednesday\" ---\n$dateOld = \"2025-11-18 Tuesday\"\n$dateNew = \"2025-11-19 Wednesday\"\n\nforeach ($para in $d.Paragraphs) {\n    $r = $para.Range\n    if ($r.Text.TrimEnd(\"`r\", \"`a\") -eq $dateOld) {\n        $r.Text = $dateNew\n        break\n    }\n}\n\n# --- 2) Table answer cells ---\n$t = $d.Tables.Item(1)\n\n# (Row, Col) are 1-based COM indices; only the 5 \"data\" rows (1, 5, 9, 13, 17)\n# of the 20-row table carry text, the rest are blank spacer rows.\n$cellUpdates = @(\n    @{ Row = 1; Col = 1; NewText = \"92\u00f72=46, 0\" }\n    @{ Row = 1; Col = 2; NewText = \"23\u00f73=7, 2\" }\n    @{ Row = 1; Col = 3; NewText = \"98\u00f72=49, 0\" }\n    @{ Row = 1; Col = 4; NewText = \"33\u00f72=16, 1\" }\n    @{ Row = 1; Col = 5; NewText = \"23\u00f75=4, 3\" }\n    @{ Row = 5; Col = 1; NewText = \"48\u00f76=8, 0\" }\n    @{ Row = 5; Col = 2; NewText = \"50\u00f74=12, 2\" }\n    @{ Row = 5; Col = 3; NewText = \"97\u00f72=48, 1\" }\n    @{ Row = 5; Col = 4; NewText = \"91\u00f77=13, 0\" }\n    @{ Row = 5; Col = 5; NewText = \"90\u00f75=18, 0\" }\n    @{ Row = 9; Col = 1; NewText = \"81\u00f79=9, 0\" }\n    @{ Row = 9; Col = 2; NewText = \"65\u00f79=7, 2\" }\n    @{ Row = 9; Col = 3; NewText = \"89\u00f73=29, 2\" }\n    @{ Row = 9; Col = 4; NewText = \"73\u00f79=8, 1\" }\n    @{ Row = 9; Col = 5; NewText = \"86\u00f77=12, 2\" }\n    @{ Row = 13; Col = 1; NewText = \"15\u00f77=2, 1\" }\n    @{ Row = 13; Col = 2; NewText = \"91\u00f73=30, 1\" }\n    @{ Row = 13; Col = 3; NewText = \"43\u00f76=7, 1\" }\n    @{ Row = 13; Col = 4; NewText = \"93\u00f77=13, 2\" }\n    @{ Row = 13; Col = 5; NewText = \"64\u00f79=7, 1\" }\n    @{ Row = 17; Col = 1; NewText = \"56\u00f73=18, 2\" }\n    @{ Row = 17; Col = 2; NewText = \"44\u00f78=5, 4\" }\n    @{ Row = 17; Col = 3; NewText = \"65\u00f72=32, 1\" }\n    @{ Row = 17; Col = 4; NewText = \"31\u00f73=10, 1\" }\n    @{ Row = 17; Col = 5; NewText = \"14\u00f79=1, 5\" }\n)\n\nforeach ($update in $cellUpdates) {\n    $t.Cell($update.Row, $update.Col).Range.Text = $update.NewText\n}\n"}
